# customer_template.xlsx — add "companiesPrimaryNames" column (M) to the
# header row, give it its own (Arial 11) font/style, widen the new column,
# and move the active selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell M1 (new shared string "companiesPrimaryNames"), styled
# with a fresh Arial/11/black font (mirrors K1/L1's existing font rules,
# just swapping Calibri for Arial).
$ws.Range("M1").Value = "companiesPrimaryNames"
$ws.Range("M1").Font.Name = "Arial"
$ws.Range("M1").Font.Size = 11
$ws.Range("M1").Font.Color = 0

# Column M gets its own width (~20.81 characters); this also splits the
# former 11-25 "K:Y" width-10.56 column band around it.
$ws.Columns.Item(13).ColumnWidth = 20.1

# Restore the author's view state: scrolled so column F is left-most,
# with H10 the active cell.
$ws.Range("H10").Select()
